# Commit: "Add pagination; change image URL property to dropbox file path"
#
# 1) Rename the "Image URL" column header (P1) on the Designs sheet to
#    "Dropbox Image Path".
# 2) Replace the Dropbox *share-link URLs* in column P (rows 2-9) with the
#    internal Dropbox *file path* for each design image
#    ("/IP Design Library/IP New Designs_2023 Spring-Summer/<id>.jpg"),
#    using a new, unique path per row (1645.jpg ... 1652.jpg) instead of the
#    original 3 repeated URLs.
# 3) Update the sheet's view/selection ("pagination") - the active cell
#    moves from P15 to I21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Designs")
$ws.Activate()

# --- Header -----------------------------------------------------------
$ws.Range("P1").Value = "Dropbox Image Path"

# --- Body: one unique Dropbox file path per row ------------------------
$paths = @(
    "/IP Design Library/IP New Designs_2023 Spring-Summer/1645.jpg",
    "/IP Design Library/IP New Designs_2023 Spring-Summer/1646.jpg",
    "/IP Design Library/IP New Designs_2023 Spring-Summer/1647.jpg",
    "/IP Design Library/IP New Designs_2023 Spring-Summer/1648.jpg",
    "/IP Design Library/IP New Designs_2023 Spring-Summer/1649.jpg",
    "/IP Design Library/IP New Designs_2023 Spring-Summer/1650.jpg",
    "/IP Design Library/IP New Designs_2023 Spring-Summer/1651.jpg",
    "/IP Design Library/IP New Designs_2023 Spring-Summer/1652.jpg"
)

for ($i = 0; $i -lt $paths.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 16).Value = $paths[$i]
}

# --- View / selection ("pagination") -----------------------------------
$ws.Range("I21").Select()
